$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 272 (existing rows 272-299 shift down
# to 274-301), matching the target sheet's dimension A1:T301.
$ws.Rows("272:273").Insert()

# Populate the two newly inserted rows with the new "Toscana" variety
# records (same Mercado/Region/Producto context as every other row).
$ws.Range("A272").Value = 7
$ws.Range("B272").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C272").Value = "Ñuble"
$ws.Range("D272").Value = 44918
$ws.Range("E272").Value = 16
$ws.Range("F272").Value = "Fruta"
$ws.Range("G272").Value = 100103
$ws.Range("H272").Value = "Frutos de hueso (carozo)"
$ws.Range("I272").Value = 100103004
$ws.Range("J272").Value = "Durazno"
$ws.Range("K272").Value = "Toscana"
$ws.Range("L272").Value = "Primera"
$ws.Range("M272").Value = 160
$ws.Range("N272").Value = 16000
$ws.Range("O272").Value = 17000
$ws.Range("P272").Value = 16500
$ws.Range("Q272").Value = "$/caja 15 kilos granel"
$ws.Range("R272").Value = "Región de O'Higgins"
$ws.Range("S272").Value = 1100
$ws.Range("T272").Value = 15

$ws.Range("A273").Value = 7
$ws.Range("B273").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C273").Value = "Ñuble"
$ws.Range("D273").Value = 44918
$ws.Range("E273").Value = 16
$ws.Range("F273").Value = "Fruta"
$ws.Range("G273").Value = 100103
$ws.Range("H273").Value = "Frutos de hueso (carozo)"
$ws.Range("I273").Value = 100103004
$ws.Range("J273").Value = "Durazno"
$ws.Range("K273").Value = "Toscana"
$ws.Range("L273").Value = "Segunda"
$ws.Range("M273").Value = 80
$ws.Range("N273").Value = 15000
$ws.Range("O273").Value = 15000
$ws.Range("P273").Value = 15000
$ws.Range("Q273").Value = "$/caja 15 kilos granel"
$ws.Range("R273").Value = "Región de O'Higgins"
$ws.Range("S273").Value = 1000
$ws.Range("T273").Value = 15
